$d = $word.ActiveDocument

$RFONTS = '<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>'
$ENDASH = [char]0x2013

function Wrap-Body([string]$innerXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
           $innerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

function R([string]$text, [bool]$preserve = $false) {
    if ($preserve) {
        return '<w:r><w:rPr>' + $RFONTS + '</w:rPr><w:t xml:space="preserve">' + $text + '</w:t></w:r>'
    } else {
        return '<w:r><w:rPr>' + $RFONTS + '</w:rPr><w:t>' + $text + '</w:t></w:r>'
    }
}

function SpellWord([string]$text) {
    return '<w:proofErr w:type="spellStart"/>' + (R $text) + '<w:proofErr w:type="spellEnd"/>'
}

# ---------------------------------------------------------------------------
# Paragraph 4: "Fixed selectivity to allow for fixing (or not) the maximum selex to 1"
#   -> split the second run so "selex" carries spell-check proofErr markers.
# ---------------------------------------------------------------------------
$p4pPr = '<w:pPr><w:pStyle w:val="NoSpacing"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr>' + $RFONTS + '</w:rPr></w:pPr>'
$p4 = '<w:p>' + $p4pPr + `
        (R 'F') + `
        (R 'ixed selectivity to allow for fixing (or not) the maximum ' $true) + `
        (SpellWord 'selex') + `
        (R ' to 1' $true) + `
      '</w:p>'
$d.Paragraphs(4).Range.InsertXML((Wrap-Body $p4))

# ---------------------------------------------------------------------------
# Paragraph 6: "Corrected the Gmacs_in.Ctl AND protected MLAState"
#   -> split so "Gmacs_in.Ctl" and "MLAState" each carry spell-check proofErr markers.
# ---------------------------------------------------------------------------
$p6pPr = '<w:pPr><w:pStyle w:val="NoSpacing"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr>' + $RFONTS + '</w:rPr></w:pPr>'
$p6 = '<w:p>' + $p6pPr + `
        (R 'Corrected the ' $true) + `
        (SpellWord 'Gmacs_in.Ctl') + `
        (R ' AND protected ' $true) + `
        (SpellWord 'MLAState') + `
      '</w:p>'
$d.Paragraphs(6).Range.InsertXML((Wrap-Body $p6))

# ---------------------------------------------------------------------------
# Paragraph 14: "Bug-fix - correct a condition on fhit in clalc_brute_equilibrium() and
# tempZ1 used in calc_predicted_project()" -> split with proofErr markers, then add two
# new bullet paragraphs describing GMACS_V_2_01_J, then (together with the trailing blank
# paragraph at the end of the document) a final paragraph stating that GMACS_V_2_01_J is
# the basis of the unified code.
# ---------------------------------------------------------------------------
$p14pPr = '<w:pPr><w:pStyle w:val="NoSpacing"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr>' + $RFONTS + '</w:rPr></w:pPr>'
$p14 = '<w:p>' + $p14pPr + `
        (R ('Bug-fix ' + $ENDASH + ' correct a condition on ') $true) + `
        (SpellWord 'fhit') + `
        (R ' in ' $true) + `
        (SpellWord 'clalc_brute_equilibrium') + `
        (R '() and tempZ1 used in ' $true) + `
        (SpellWord 'calc_predicted_project') + `
        (R '()') + `
      '</w:p>'

$p14bpPr = '<w:pPr><w:pStyle w:val="NoSpacing"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr>' + $RFONTS + '</w:rPr></w:pPr>'
$p14b = '<w:p>' + $p14bpPr + `
        (R 'GMACS_V_2_01_') + `
        (R 'J: ' $true) + `
        (R 'GMACS_V_2_01_') + `
        (R 'I  + ' $true) + `
      '</w:p>'

$p14cpPr = '<w:pPr><w:pStyle w:val="NoSpacing"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr>' + $RFONTS + '</w:rPr></w:pPr>'
$p14c = '<w:p>' + $p14cpPr + `
        (R 'Option to consider terminally molting life history') + `
      '</w:p>'

# original trailing blank paragraph (unchanged formatting)
$p15pPr = '<w:pPr><w:pStyle w:val="NoSpacing"/><w:ind w:left="720"/><w:rPr>' + $RFONTS + '</w:rPr></w:pPr>'
$p15 = '<w:p>' + $p15pPr + '</w:p>'

$p16pPr = '<w:pPr><w:pStyle w:val="NoSpacing"/><w:rPr>' + $RFONTS + '</w:rPr></w:pPr>'
$p16 = '<w:p>' + $p16pPr + `
        (R 'GMACS_V_2_01_') + `
        (R 'J') + `
        (R '  = ' $true) + `
        (R 'UNIFIED CODE') + `
      '</w:p>'

# Replace paragraph 14 through the last (trailing blank) paragraph in one shot, so the
# replacement range extends to the true end of the document body and no stray blank
# paragraph is left behind.
$startRange = $d.Paragraphs(14).Range.Start
$endRange = $d.Paragraphs($d.Paragraphs.Count).Range.End
$tailRange = $d.Range($startRange, $endRange)
$tailRange.InsertXML((Wrap-Body ($p14 + $p14b + $p14c + $p15 + $p16)))

Write-Host "Done. Final paragraph count: " $d.Paragraphs.Count
